$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.351.04'
$ws.Range('D3').Value = '2.222.89'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = "'110.82"
$ws.Range('E5').Value = '  -8.19%  '
$ws.Range('D6').Value = "'289.64"
$ws.Range('E6').Value = '  +7.74%  '
$ws.Range('E7').Value = '  -2.61%  '
$ws.Range('E8').Value = '  -0.39%  '
$ws.Range('D9').Value = "'0.598"
$ws.Range('E9').Value = '  -3.63%  '
$ws.Range('D10').Value = "'43.50"
$ws.Range('E10').Value = '  -8.43%  '
$ws.Range('D11').Value = "'0.0908"
$ws.Range('E11').Value = '  -3.57%  '
$ws.Range('D12').Value = "'54.15"
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').Value = "'8.60"
$ws.Range('E13').Value = '  -8.45%  '
$ws.Range('E14').Value = '  +9.95%  '
$ws.Range('E15').Value = '  -2.97%  '
$ws.Range('D16').Value = "'14.87"
$ws.Range('E16').Value = '  -5.81%  '
$ws.Range('D17').Value = '2.559.56'
$ws.Range('E17').Value = '  -1.91%  '
$ws.Range('D18').Value = '2.231.11'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('D19').Value = '42.347.16'
$ws.Range('E19').Value = '  -2.84%  '
$ws.Range('D20').Value = "'7.11"
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('D21').Value = "'0.0000104"
$ws.Range('E21').Value = '  -4.51%  '
$ws.Range('D22').Value = "'72.68"
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').Value = "'3.34"
$ws.Range('E23').Value = '  +13.05%  '
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('D25').Value = "'230.19"
$ws.Range('E25').Value = '  -2.01%  '
$ws.Range('E26').Value = '  -7.04%  '
$ws.Range('D28').Value = "'11.40"
$ws.Range('E28').Value = '  -6.98%  '
$ws.Range('D29').Value = "'2.20"
$ws.Range('E29').Value = '  -2.52%  '
$ws.Range('D30').Value = "'37.57"
$ws.Range('E30').Value = '  -10.92%  '
$ws.Range('D31').Value = "'173.44"
$ws.Range('E32').Value = '  -7.63%  '
$ws.Range('D33').Value = "'20.83"
$ws.Range('E33').Value = '  -3.20%  '
$ws.Range('E34').Value = '  -4.16%  '
$ws.Range('E35').Value = '  -1.99%  '
$ws.Range('D36').Value = "'4.95"
$ws.Range('E36').Value = '  +5.06%  '
$ws.Range('D37').Value = "'4.22"
$ws.Range('E37').Value = '  -5.92%  '
$ws.Range('E38').Value = '  -3.57%  '
$ws.Range('D39').Value = "'0.0371"
$ws.Range('E39').Value = '  -2.39%  '
$ws.Range('E40').Value = '  -4.70%  '
$ws.Range('D41').Value = "'73.63"
$ws.Range('E41').Value = '  +1.70%  '
$ws.Range('D42').Value = "'2.38"
$ws.Range('E42').Value = '  -6.52%  '
$ws.Range('D43').Value = "'0.230"
$ws.Range('E43').Value = '  -4.81%  '
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('E45').Value = '  -10.89%  '
$ws.Range('E46').Value = '  -4.50%  '
$ws.Range('E47').Value = '  -6.85%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').Value = "'1.27"
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = "'1.67"
$ws.Range('E49').Value = '  +4.71%  '
$ws.Range('D50').Value = "'101.33"
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('E51').Value = '  -2.00%  '
